$wb = $excel.ActiveWorkbook

# Column-width in Excel's ColumnWidth units is stored-width minus 5/6 (0.8333...);
# the target stored width is exactly 40 characters for columns I, J and P.
$targetColumnWidth = 40 - (5/6)

$sheetInfo = @(
    @{
        Name = "zh-cn"
        HandbackDateTime = "2016-11-07 06:54:23"
        TargetFile = "a14b8ca5-f559-4148-9701-350adb07cd9e.md"
        HandbackFile = "a14b8ca5-f559-4148-9701-350adb07cd9e.24e393409de0d4860379416983f101ce9b2505f3.zh-cn.xlf"
    },
    @{
        Name = "de-de"
        HandbackDateTime = "2016-11-07 06:54:42"
        TargetFile = "a14b8ca5-f559-4148-9701-350adb07cd9e.md"
        HandbackFile = "a14b8ca5-f559-4148-9701-350adb07cd9e.24e393409de0d4860379416983f101ce9b2505f3.de-de.xlf"
    }
)

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dab2c4adb2271edb501907d3d10923f49d4f0af6/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7100a5fa181df88894f3deac028304d24a80d08b/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md."
$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7100a5fa181df88894f3deac028304d24a80d08b/e2e/a14b8ca5-f559-4148-9701-350adb07cd9e.md"

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen "Latest Target File" (I), "Latest Handback File" (J) and "Error Detail" (P) columns.
    $ws.Columns.Item(9).ColumnWidth = $targetColumnWidth
    $ws.Columns.Item(10).ColumnWidth = $targetColumnWidth
    $ws.Columns.Item(16).ColumnWidth = $targetColumnWidth

    # Fill in the handback-report values for row 5 (the a14b8ca5 file).
    $ws.Range("J5").Value = $info.HandbackFile
    $ws.Range("K5").Value = $info.HandbackDateTime
    $ws.Range("P5").Value = $errorDetail

    # "Latest Target File" becomes a hyperlink to the latest version of the source md file.
    $ws.Hyperlinks.Add($ws.Range("I5"), $targetMdUrl, [System.Type]::Missing, [System.Type]::Missing, $info.TargetFile)
    $ws.Range("I5").Font.Underline = $true
    $ws.Range("I5").Font.Color = 15570276
}
